$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.723.00'
$ws.Range("E2").Value = '  +1.87%  '

$ws.Range("D3").Value = '1.900.52'
$ws.Range("E3").Value = '  +2.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("E7").Value = '  +1.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2844'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06556'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.97%  '

$ws.Range("D10").Value = '1.897.91'
$ws.Range("E10").Value = '  +2.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07478'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.17%  '

$ws.Range("E12").Value = '  +2.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.118'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6688'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.18%  '

$ws.Range("D16").Value = '30.689.12'
$ws.Range("E16").Value = '  +1.94%  '

$ws.Range("E17").Value = '  +0.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007623'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.75%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.142.48'
$ws.Range("E21").Value = '  +1.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.309'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.45%  '

$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.237'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.283'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.45%  '

$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.956'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.81%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.400'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1007'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.360'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.034'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.01%  '

$ws.Range("E33").Value = '  +1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.223'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7541'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.99%  '

$ws.Range("E36").Value = '  +0.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01886'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.657'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9179'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.087'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.74%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4301'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.62%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.831'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.005'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.423'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("E47").Value = '  -2.48%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.487'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.25%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.023'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05663'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '
